$p = $ppt.ActivePresentation

# The Handout Master and Notes Master each carry an auto-updating
# "datetime1" date field whose last-computed value (10/5/15) was cached
# in the OOXML. Re-point that cached date to 12/21/15, as happens when
# the deck is reopened/resaved on a later date.
$p.HandoutMaster.HeadersFooters.DateAndTime.Text = "12/21/15"
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "12/21/15"
